# Cover-letter text touch-ups: add back the detail that an earlier model
# was developed, note the speed comparison to RVD2, and add the missing
# article "the" before "NCBI Sequence Read Archive".

$d = $word.ActiveDocument

# 1) "We have previously developed and published" ->
#    "We have previously developed a variant detection model and published"
$d.Content.Find.Execute(
    "developed and published",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "developed a variant detection model and published",
    2
)

# 2) "...more computationally efficient on tests of low coverage..." ->
#    "...more computationally efficient than RVD2 on tests of low coverage..."
$d.Content.Find.Execute(
    "more computationally efficient on tests",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "more computationally efficient than RVD2 on tests",
    2
)

# 3) "...is available on NCBI Sequence Read Archive." ->
#    "...is available on the NCBI Sequence Read Archive."
$d.Content.Find.Execute(
    "available on NCBI Sequence Read Archive",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "available on the NCBI Sequence Read Archive",
    2
)
